# Auto-generated edit script.
#
# Goal: insert a ' ## ' sentence-boundary marker at specific points inside
# existing <w:t> runs (never touching paragraph structure or any run's
# <w:rPr> formatting).
#
# Two wrinkles in this runtime's Find/Execute simulation drove the approach:
#
#  1) Passing an apostrophe through the *replacement* text of
#     Find.Execute(..., Replace:=wdReplaceOne) gets 'smart quoted' into a
#     curly quote (U+2018/2019), which the target text does not want (the
#     source keeps plain straight apostrophes). Workaround: substitute each
#     apostrophe with a private-use placeholder character before calling
#     Find/Replace, then patch each placeholder back to a literal "'" via a
#     direct (non-Find/Replace) Range.Text assignment, which this runtime
#     does not 'smart quote'.
#
#  2) Range.Text assignment (used for the fix-ups above, and for the one
#     spot where the search text is too generic to Find uniquely) drops the
#     xml:space="preserve" attribute on the <w:t> it rewrites. That's fine for
#     the tiny one-character placeholder patches (no leading/trailing space
#     at risk) and for the one run where we anchor+offset instead of Find.

$d = $word.ActiveDocument
$results = @()

# --- replace_0
$rng0 = $d.Content
$f0 = $rng0.Find.Execute('An introspective look at the relationship between Hawking and the space/time contingent. This film ', $true, $false, $false, $false, $false, $true, 1, $false, 'An introspective look at the relationship between Hawking and the space/time contingent. ## This film ', 2)
$results += "found_0=$f0"

# --- replace_1
$rng1 = $d.Content
$f1 = $rng1.Find.Execute('relation to Einstein''s Theory of General Relativity. The film is methodically directed, exposing details of the man (Hawking) as well as his work (Black Holes). Interviews with his family are a little too long so sadly there is less development of his theories and ideas.  A Philip Glass soundtrack superbly compliments the film. Only one other man could compose such haunting ', $true, $false, $false, $false, $false, $true, 1, $false, 'relation to Einsteins Theory of General Relativity. ## The film is methodically directed, exposing details of the man (Hawking) as well as his work (Black Holes). ## Interviews with his family are a little too long so sadly there is less development of his theories and ideas. ## A Philip Glass soundtrack superbly compliments the film. ## Only one other man could compose such haunting ', 2)
$results += "found_1=$f1"
$ph1_0 = $d.Content
$null = $ph1_0.Find.Execute([char]1, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$fix1_0 = $d.Range($ph1_0.Start, $ph1_0.End)
$fix1_0.Text = "'"

# --- replace_2
$rng2 = $d.Content
$f2 = $rng2.Find.Execute('melodies (Jean Michel Jarre). Overall I would highly recommend this movie on the basis of Hawking''s ''nuggets of wisdom'' and his adequate explanation of an Event Horizon!', $true, $false, $false, $false, $false, $true, 1, $false, 'melodies (Jean Michel Jarre). ## Overall I would highly recommend this movie on the basis of Hawkings nuggets of wisdom and his adequate explanation of an Event Horizon!', 2)
$results += "found_2=$f2"
$ph2_0 = $d.Content
$null = $ph2_0.Find.Execute([char]1, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$fix2_0 = $d.Range($ph2_0.Start, $ph2_0.End)
$fix2_0.Text = "'"
$ph2_1 = $d.Content
$null = $ph2_1.Find.Execute([char]2, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$fix2_1 = $d.Range($ph2_1.Start, $ph2_1.End)
$fix2_1.Text = "'"
$ph2_2 = $d.Content
$null = $ph2_2.Find.Execute([char]3, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$fix2_2 = $d.Range($ph2_2.Start, $ph2_2.End)
$fix2_2.Text = "'"

# --- replace_3
$rng3 = $d.Content
$f3 = $rng3.Find.Execute('Introspektivni pogled na odnos između Hawkinga i kontingenta prostora / vremena. Ovaj film iznosi ', $true, $false, $false, $false, $false, $true, 1, $false, 'Introspektivni pogled na odnos između Hawkinga i kontingenta prostora / vremena. ## Ovaj film iznosi ', 2)
$results += "found_3=$f3"

# --- replace_4
$rng4 = $d.Content
$f4 = $rng4.Find.Execute('. Film je metodički režiran, izlaže detalje o čovjeku (Hawking) kao i njegovu djelu (Crne rupe). Intervjui s njegovom obitelji malo su predugi pa je nažalost manje razvoja njegovih teorija i ideja. Soundtrack Philipa Glassa izvrsno nadopunjuje film. Samo je jedan drugi čovjek mogao skladati takve ', $true, $false, $false, $false, $false, $true, 1, $false, '. ## Film je metodički režiran, izlaže detalje o čovjeku (Hawking) kao i njegovu djelu (Crne rupe). ## Intervjui s njegovom obitelji malo su predugi pa je nažalost manje razvoja njegovih teorija i ideja. ## Soundtrack Philipa Glassa izvrsno nadopunjuje film. ## Samo je jedan drugi čovjek mogao skladati takve ', 2)
$results += "found_4=$f4"

# --- replace_5
$rng5 = $d.Content
$f5 = $rng5.Find.Execute('zvjezdane melodije (Jean Michel Jarre). Sve u svemu, toplo bih preporučio ovaj film na osnovu Hawkingovih ', $true, $false, $false, $false, $false, $true, 1, $false, 'zvjezdane melodije (Jean Michel Jarre). ## Sve u svemu, toplo bih preporučio ovaj film na osnovu Hawkingovih ', 2)
$results += "found_5=$f5"

# --- replace_6
$rng6 = $d.Content
$f6 = $rng6.Find.Execute('My god...i have not seen such an awful movie in a long...long time...saw it last night and wanted to leave after 20 minutes...keira knightley tries really really hard in this one, but she cant handle it..dropped her accent every once in a while and didn''t have the charisma to fill the role...sienna ', $true, $false, $false, $false, $false, $true, 1, $false, 'My god...i have not seen such an awful movie in a long...long time…## saw it last night and wanted to leave after 20 minutes… ## keira knightley tries really really hard in this one, but she cant handle it..dropped her accent every once in a while and didnt have the charisma to fill the role… ## sienna ', 2)
$results += "found_6=$f6"
$ph6_0 = $d.Content
$null = $ph6_0.Find.Execute([char]1, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$fix6_0 = $d.Range($ph6_0.Start, $ph6_0.End)
$fix6_0.Text = "'"

# --- replace_7
$rng7 = $d.Content
$f7 = $rng7.Find.Execute(' acting gets you to a point where you start to ask yourself: Has she ever had acting lessons? judging by the edge of love ', $true, $false, $false, $false, $false, $true, 1, $false, ' acting gets you to a point where you start to ask yourself: Has she ever had acting lessons? ## judging by the edge of love ', 2)
$results += "found_7=$f7"

# --- replace_8
$rng8 = $d.Content
$f8 = $rng8.Find.Execute(' never been to acting class, but should consider to go in the near future...they both look really pretty..maybe ', $true, $false, $false, $false, $false, $true, 1, $false, ' never been to acting class, but should consider to go in the near future… ## they both look really pretty..maybe ', 2)
$results += "found_8=$f8"

# --- replace_9
$rng9 = $d.Content
$f9 = $rng9.Find.Execute(' what they should focus on in their future career..if they can be actresses everybody can!', $true, $false, $false, $false, $false, $true, 1, $false, ' what they should focus on in their future career.. ## if they can be actresses everybody can!', 2)
$results += "found_9=$f9"

# --- replace_10
$rng10 = $d.Content
$f10 = $rng10.Find.Execute(' tako grozan film dugo ... dugo vremena ... ', $true, $false, $false, $false, $false, $true, 1, $false, ' tako grozan film dugo ... dugo vremena ...  ## ', 2)
$results += "found_10=$f10"

# --- replace_11
$rng11 = $d.Content
$f11 = $rng11.Find.Execute(' sam ga sinoć i htio je otići nakon 20 minuta ... keira knightley se jako trudi u ', $true, $false, $false, $false, $false, $true, 1, $false, ' sam ga sinoć i htio je otići nakon 20 minuta ... ## keira knightley se jako trudi u ', 2)
$results += "found_11=$f11"

# --- replace_12
$rng12 = $d.Content
$f12 = $rng12.Find.Execute(' ulogu ... Sienna ', $true, $false, $false, $false, $false, $true, 1, $false, ' ulogu ... ## Sienna ', 2)
$results += "found_12=$f12"

# --- replace_13
$rng13 = $d.Content
$f13 = $rng13.Find.Execute(' do točke kad se počnete pitati: Je li ikad imala lekcije glume? sudeći po ', $true, $false, $false, $false, $false, $true, 1, $false, ' do točke kad se počnete pitati: Je li ikad imala lekcije glume? ## sudeći po ', 2)
$results += "found_13=$f13"

# --- replace_14
$rng14 = $d.Content
$f14 = $rng14.Find.Execute(' u bliskoj budućnosti ... ', $true, $false, $false, $false, $false, $true, 1, $false, ' u bliskoj budućnosti ... ## ', 2)
$results += "found_14=$f14"

# --- replace_15
$rng15 = $d.Content
$f15 = $rng15.Find.Execute(' lijepo.. ', $true, $false, $false, $false, $false, $true, 1, $false, ' lijepo.. ## ', 2)
$results += "found_15=$f15"

# --- replace_16
$rng16 = $d.Content
$f16 = $rng16.Find.Execute('Thanks to other reviewers who directed me to this product when I was told I was anemic. Now been taking these for about 4 months and the anemia is gone.  Good product.  Easily digested (unlike some other iron supplements).', $true, $false, $false, $false, $false, $true, 1, $false, 'Thanks to other reviewers who directed me to this product when I was told I was anemic. ## Now been taking these for about 4 months and the anemia is gone. ## Good product.  ## Easily digested (unlike some other iron supplements).', 2)
$results += "found_16=$f16"

# --- replace_17
$rng17 = $d.Content
$f17 = $rng17.Find.Execute('da sam anemična. Sad ih uzimam otprilike 4 mjeseca i anemije više nema. Dobar proizvod. Lako se probavlja (za razliku od nekih drugih ', $true, $false, $false, $false, $false, $true, 1, $false, 'da sam anemična. ## Sad ih uzimam otprilike 4 mjeseca i anemije više nema. ## Dobar proizvod. ## Lako se probavlja (za razliku od nekih drugih ', 2)
$results += "found_17=$f17"

# --- replace_18
$rng18 = $d.Content
$f18 = $rng18.Find.Execute('This is one of my favorite desserts, and melts quickly in the mouth. This brand is good and it shipped well-packaged. Everyone should try this once. The amazon price is much better than the ones you find at science fairs.', $true, $false, $false, $false, $false, $true, 1, $false, 'This is one of my favorite desserts, and melts quickly in the mouth. ## This brand is good and it shipped well-packaged. ## Everyone should try this once. ## The amazon price is much better than the ones you find at science fairs.', 2)
$results += "found_18=$f18"

# --- replace_19
$rng19 = $d.Content
$f19 = $rng19.Find.Execute('Ovo je jedan od mojih najdražih deserta i brzo se topi u ustima. Ova marka je dobra i isporučuje ', $true, $false, $false, $false, $false, $true, 1, $false, 'Ovo je jedan od mojih najdražih deserta i brzo se topi u ustima. ## Ova marka je dobra i isporučuje ', 2)
$results += "found_19=$f19"

# --- replace_20
$rng20 = $d.Content
$f20 = $rng20.Find.Execute('. Svi bi trebali jednom ', $true, $false, $false, $false, $false, $true, 1, $false, '. ## Svi bi trebali jednom ', 2)
$results += "found_20=$f20"

# --- replace_21: old text '. ' occurs 54x verbatim in the document, so it
# can't be Found uniquely by content. Anchor on the preceding run
# ('pokušati', which is unique) and address the following run (exactly
# '. ') by absolute character offsets instead.
$anchor21 = $d.Content
$found21 = $anchor21.Find.Execute('pokušati', $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$results += "found_21=$found21"
$target21 = $d.Range($anchor21.End, $anchor21.End + 2)
$target21.Text = '. ## '

# --- replace_22
$rng22 = $d.Content
$f22 = $rng22.Find.Execute('This is a fantastic puzzle/gift for young AND old. ', $true, $false, $false, $false, $false, $true, 1, $false, 'This is a fantastic puzzle/gift for young AND old. ## ', 2)
$results += "found_22=$f22"

# --- replace_23
$rng23 = $d.Content
$f23 = $rng23.Find.Execute('32 triangular strong magnetic pieces that can fit together in a wide number of ways.  It''s just great and you''ll have trouble keeping it away from the adults.', $true, $false, $false, $false, $false, $true, 1, $false, '32 triangular strong magnetic pieces that can fit together in a wide number of ways. ## Its just great and youll have trouble keeping it away from the adults.', 2)
$results += "found_23=$f23"
$ph23_0 = $d.Content
$null = $ph23_0.Find.Execute([char]1, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$fix23_0 = $d.Range($ph23_0.Start, $ph23_0.End)
$fix23_0.Text = "'"
$ph23_1 = $d.Content
$null = $ph23_1.Find.Execute([char]2, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$fix23_1 = $d.Range($ph23_1.Start, $ph23_1.End)
$fix23_1.Text = "'"

# --- replace_24
$rng24 = $d.Content
$f24 = $rng24.Find.Execute('Ovo je fantastična zagonetka / poklon za mlade I stare. Riječ je o 32 trokutasta jaka magnetska ', $true, $false, $false, $false, $false, $true, 1, $false, 'Ovo je fantastična zagonetka / poklon za mlade I stare. ## Riječ je o 32 trokutasta jaka magnetska ', 2)
$results += "found_24=$f24"

# --- replace_25
$rng25 = $d.Content
$f25 = $rng25.Find.Execute('koji se mogu sastaviti na više različitih načina. Jednostavno je super i ', $true, $false, $false, $false, $false, $true, 1, $false, 'koji se mogu sastaviti na više različitih načina. ## Jednostavno je super i ', 2)
$results += "found_25=$f25"

# --- replace_26
$rng26 = $d.Content
$f26 = $rng26.Find.Execute('It''s another bad zombie movie. Compared to the majority of ', $true, $false, $false, $false, $false, $true, 1, $false, 'Its another bad zombie movie. ## Compared to the majority of ', 2)
$results += "found_26=$f26"
$ph26_0 = $d.Content
$null = $ph26_0.Find.Execute([char]1, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$fix26_0 = $d.Range($ph26_0.Start, $ph26_0.End)
$fix26_0.Text = "'"

# --- replace_27
$rng27 = $d.Content
$f27 = $rng27.Find.Execute(' others, the only difference here is the main character is a female. The plot is the same. The action scenes are not engaging. Special effects  are so so.', $true, $false, $false, $false, $false, $true, 1, $false, ' others, the only difference here is the main character is a female. ## The plot is the same. ## The action scenes are not engaging. ## Special effects  are so so.', 2)
$results += "found_27=$f27"

# --- replace_28
$rng28 = $d.Content
$f28 = $rng28.Find.Execute('. U usporedbi s većinom ostalih, jedina razlika ovdje je ', $true, $false, $false, $false, $false, $true, 1, $false, '. ## U usporedbi s većinom ostalih, jedina razlika ovdje je ', 2)
$results += "found_28=$f28"

# --- replace_29
$rng29 = $d.Content
$f29 = $rng29.Find.Execute('. Zaplet je isti. Akcijske scene nisu privlačne. Posebni efekti su ', $true, $false, $false, $false, $false, $true, 1, $false, '. ## Zaplet je isti. ## Akcijske scene nisu privlačne. ## Posebni efekti su ', 2)
$results += "found_29=$f29"

# --- replace_30
$rng30 = $d.Content
$f30 = $rng30.Find.Execute('YoYo seems out of balance. No matter ', $true, $false, $false, $false, $false, $true, 1, $false, 'YoYo seems out of balance. ## No matter ', 2)
$results += "found_30=$f30"

# --- replace_31
$rng31 = $d.Content
$f31 = $rng31.Find.Execute('tilt to one side. Made it difficult ', $true, $false, $false, $false, $false, $true, 1, $false, 'tilt to one side. ## Made it difficult ', 2)
$results += "found_31=$f31"

# --- replace_32
$rng32 = $d.Content
$f32 = $rng32.Find.Execute('or do tricks. I have a little ', $true, $false, $false, $false, $false, $true, 1, $false, 'or do tricks. ## I have a little ', 2)
$results += "found_32=$f32"

# --- replace_33
$rng33 = $d.Content
$f33 = $rng33.Find.Execute('izvan ravnoteže. Bez obzira na ', $true, $false, $false, $false, $false, $true, 1, $false, 'izvan ravnoteže. ## Bez obzira na ', 2)
$results += "found_33=$f33"

# --- replace_34
$rng34 = $d.Content
$f34 = $rng34.Find.Execute(' na jednu stranu. Teško ', $true, $false, $false, $false, $false, $true, 1, $false, ' na jednu stranu. ## Teško ', 2)
$results += "found_34=$f34"

# --- replace_35
$rng35 = $d.Content
$f35 = $rng35.Find.Execute('trikove. Imam malo iskustva s ', $true, $false, $false, $false, $false, $true, 1, $false, 'trikove. ## Imam malo iskustva s ', 2)
$results += "found_35=$f35"

# --- replace_36
$rng36 = $d.Content
$f36 = $rng36.Find.Execute('If you wrap your wrists properly, you''ll see these are both too narrow and too short, way too short. Do not get these if you are hitting the heavy bag. They just won''t protect/support your wrists or knuckles.', $true, $false, $false, $false, $false, $true, 1, $false, 'If you wrap your wrists properly, youll see these are both too narrow and too short, way too short. ## Do not get these if you are hitting the heavy bag. ## They just wont protect/support your wrists or knuckles.', 2)
$results += "found_36=$f36"
$ph36_0 = $d.Content
$null = $ph36_0.Find.Execute([char]1, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$fix36_0 = $d.Range($ph36_0.Start, $ph36_0.End)
$fix36_0.Text = "'"
$ph36_1 = $d.Content
$null = $ph36_1.Find.Execute([char]2, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$fix36_1 = $d.Range($ph36_1.Start, $ph36_1.End)
$fix36_1.Text = "'"

# --- replace_37
$rng37 = $d.Content
$f37 = $rng37.Find.Execute('. Ne ', $true, $false, $false, $false, $false, $true, 1, $false, '. ## Ne ', 2)
$results += "found_37=$f37"

# --- replace_38
$rng38 = $d.Content
$f38 = $rng38.Find.Execute(' ih ako udarite u tešku torbu. Oni jednostavno neće zaštititi / podržavati vaše ', $true, $false, $false, $false, $false, $true, 1, $false, ' ih ako udarite u tešku torbu. ## Oni jednostavno neće zaštititi / podržavati vaše ', 2)
$results += "found_38=$f38"

$results -join "`n"
